$d = $word.ActiveDocument

# --- Paragraph 1: title line - update date and paper title ---
$d.Content.Find.Execute("⚡️🚀המאמר היומי של מייק -01.11.24: ⚡️🚀", $true, $false, $false, $false, $false, $true, 1, $false, "⚡️🚀המאמר היומי של מייק -31.10.24: ⚡️🚀", 2) | Out-Null
$d.Content.Find.Execute("LLMs Are In-Context Reinforcement Learners", $true, $false, $false, $false, $false, $true, 1, $false, "Understanding Transformers via N-gram Statistics", 2) | Out-Null

# --- Paragraph 2: replace body text, then append a trailing line break ---
$d.Content.Find.Execute("אני אוהב מאמרים שמשלבים כמה שיטות של ML. אסקור היום אחד כזה המציע לשדך למידת in-context עם למידה באמצעות חיזוקים או בקצרה RL. למידת in-context היא יכולת של מודל שפה ללמוד משהו חדש מכמה דוגמאות בפרופמט ללא צורך בפיין טיון. יש לא מעט הסברים ליכולת די מפתיע זו ולפעמים יכולת זו נקראה emergent capabilities.", $true, $false, $false, $false, $false, $true, 1, $false, "מאמר די נחמד ולא רגיל מבית גוגל. המאמר מחזיר אותנו לתקופה שלא מידלנו את השפה הטבעית באמצעות מודלים סטטיסטיים עם עשרות ומאות מיליארדי פרמטרים.  פעם ניסינו להשתמש ב- n-grams כדי לשערך את ההתפלגות של המילים בטקסט. כמובן גישות כאלו לא יכולות לעבוד עבור דאטהסטים בעל עשרות טריליוני טוקנים כמו שיש לנו היום אבל אולי אפשר לקחת LLMs גדולים ולבדוק האם ניתן לקרב את חיזויהם באמצעות סטטיסטיקות על n-grams. כדי לא לסבך המאמר לא בודק את זה על למידת in-context.", 2) | Out-Null
$p2 = $d.Paragraphs.Item(2)
$r2 = $p2.Range
$r2body = $d.Range($r2.Start, $r2.End - 1)
$r2body.InsertAfter([char]11)

# --- Paragraph 3: replace body text ---
$d.Content.Find.Execute("עכשיו נשאלת השאלה: איך נוכל לבחור דוגמאת להדגמה שאנו מראים למודל שפה בפרומפט למקסום ביצועיי המודל? השאלה הזו לא מאוד טריויאלית ואין עליה כרגע תשובה חד משמעית. המחברים מציעים לגשת לבעיה זו דרך למידה עם חיזוקים (סוג של). השיטה הנאיבית היא פשוט לצבור דוגמאות עד שנגמר לנו את אורך חלון ההקשר של המודל. לכל דוגמא בהדגמה אנו שומרים בבאפר את השלישיה המכילה את הדוגמא (שאלה עצמה)ֿ, תשובת המודל ומשערך של איכות התשובה (או פשוט האם התשובה נכונה או לא). ואז באינפרנס פשוט לוקחים את הדוגמאות האלו בתור פרומפט.", $true, $false, $false, $false, $false, $true, 1, $false, "וזה בדיוק מה שהמאמר הזה (שיש לו רק מחבר אחד שזה די נדיר בימינו) עושה. הוא בודק האם ניתן לחזות את הטוקן הבא שמודל שפה מאומן חוזר באמצעות סטטיסטיקה של n-grams שבאים לפניו בטקסט. במקרה הזה n-grams בנויים לא ממילים אלא מטוקנים. דרך אגב הסטטיסטיקה של n-grams אינה חייבית לכלול את כל n הטוקנים הבאים לפני הטוקן הנחזה אלא עשויה ״להכיל חורים״(כלומר יכולה לקחת טוקן i-1, i-2 i ו- i-4 עבור 3-gram - נצטרך למצע מעל טוקן i-3 בשביל כך).", 2) | Out-Null

# --- Paragraph 4: replace body text ---
$d.Content.Find.Execute("לטענת המחברים הגישה הנאיבית הזו לא עובדת משתי סיבות עיקריות. קודם כל שילוב מתמשך של אותם הפרומפטים לדוגמאות שונות מוביל לשונות גדולה בפלט של LLM (לפי המחקרים הקודמים עלולה להוביל לביצועים ירודים). הסיבה השניה טמונה בכך ששלישיות (שאלה, תשובה, לא נכון) מסבכות את המודל ולא מספקות לו מספיק מידע על איך היה צריך לענות נכון (ד״א בלמידה ניגודית יש בעיה דומה המצריכה כמות מאוד גדולה של דוגמאות שליליות בכל באץ' - כתבתי על זה לא מעט בסקירותיי).", $true, $false, $false, $false, $false, $true, 1, $false, "המחבר מצא כמה דברים מעניינים. ניתן לשערך את החיזוי של מודל שפה עם 7-gram (עבור דאטהסטים שהם בחרו) בלא מעט מקרים. בנסוף נמצא כי לטוקנים בעל שונות נמוכה (של ההתפלגות שלהם) n-grams מצליחים יותר מאשר לטוקנים בעל שונות חיזוי גבוהה. מעניין שככל שמאמנים מודל שפה יותר יותר קשה לקרב אותה עם n-grams (צריך להגדיל את n או לא משנה מה ה-n דיוק הקירוב יורד).", 2) | Out-Null

# --- Paragraph 5: replace body text with short closing remark ---
$d.Content.Find.Execute("עקב כך המחברים הציעו להכניס קצת ״אקראיות״ לבניית הפרומפטים (המחברים קוראים לזה אפיזודה בהתאם לטרמינולוגיה של RL - כל אפיזודה מורכבת מכמה שלישיות של שאלה, תשובה, נכונות התשובה) וגם להשתמש באפיזודות שקיבלו ציון ״נכון״. לכל דוגמא הם הציע קודם לדגום באקראי מהבאפר של אפיזודות בצורה אקראית ולהשתמש לכל דוגמא במדגם שונה של אפיזודות. כאמור שומרים רק את האפיזודות שבהם המודל צדק. כך פרומפט לכל שאילתה הופך להיות לא קבוע ומכיל רק דוגמאות עם תשובות נכונות. זה נקרא Explorative ICRL במאמר.", $true, $false, $false, $false, $false, $true, 1, $false, "אהבתי…", 2) | Out-Null

# --- Delete the two paragraphs that were removed entirely ---
# (search by their distinctive leading text, delete the whole paragraph range
#  including its paragraph mark so the surrounding paragraphs merge correctly)
$targets = @(
    "כמובן ש Explorative ICRL לא יעיל חישובית כי כל פעם צריך לחשב את הפרומפט מחדש (מה שלא צריך לעשות בגישה הנאיבית אך לא עובדת). המחברים שכללו את זה עם מנגנון קאשינג המאפשר לשמור מספר קבוע של פרומפרטים (מערך של אפיזודות) ולכל אפיזודה נתונה להחליט לאלו מהם להוסיף אותה. זה מקל על העלות החישובית.",
    "מאמר חמוד למרות שמשום מה לקח לי קצת זמן להבין אותו…"
)
foreach ($target in $targets) {
    for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $target) {
            $p.Range.Delete()
            break
        }
    }
}

# --- Replace the arxiv link text ---
$d.Content.Find.Execute("https://arxiv.org/pdf/2410.05362", $true, $false, $false, $false, $false, $true, 1, $false, "https://www.arxiv.org/abs/2407.12034", 2) | Out-Null
